# Practice tasks and final revisions
# Update sheet names and B-column stim/file-name values across the five
# task-order worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO --------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996104191329.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961042153344.csv"
$ws1.Range("B4").Value = "go_stims-16509961042153344.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961042313364.csv"
$ws1.Name = "GNG_TO-16509961042313364"

# --- Sheet 2: NB_TO ----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509961056573522.csv"
$ws2.Range("B3").Value = "ZB-match_6-16509961045853484.csv"
$ws2.Range("B4").Value = "OB-16509961055853527.csv"
$ws2.Range("B5").Value = "ZB-match_8-165099610441531.csv"
$ws2.Range("B6").Value = "TB-16509961057373567.csv"
$ws2.Range("B7").Value = "TB-1650996105601349.csv"
$ws2.Range("B8").Value = "ZB-match_0-16509961045131762.csv"
$ws2.Range("B9").Value = "OB-16509961053933468.csv"
$ws2.Range("B10").Value = "OB-16509961054893525.csv"
$ws2.Name = "NB_TO-1650996105753365"

# --- Sheet 3: RS_TO ------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"
$ws3.Name = "RS_TO-1650996105753365"

# --- Sheet 4: TOL_TO -----------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650996105769371.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996105753365.csv"
$ws4.Range("B4").Value = "MM_stims-16509961057853591.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996105769371.csv"
$ws4.Range("B6").Value = "MM_stims-1650996105801384.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961057853591.csv"
$ws4.Name = "TOL_TO-1650996105801384"

# --- Sheet 5: vSAT_TO ------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509961058493776.csv"
$ws5.Range("B3").Value = "SAT_stims-1650996105801384.csv"
$ws5.Range("B4").Value = "SAT_stims-1650996105817358.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961058333912.csv"
$ws5.Name = "vSAT_TO-16509961058653858"
